# Apply the "Fruta / hortaliza, semanal" weekly update:
# shift rows 16-21 up from 17-22 (for the price/volume/date/origin columns)
# and append a fresh observation into row 22.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Range("D16").Value = 44645
$ws.Range("M16").Value = 5
$ws.Range("N16").Value = 30000
$ws.Range("O16").Value = 30000
$ws.Range("P16").Value = 30000
$ws.Range("Q16").Value = '$/caja 18 kilos'
$ws.Range("R16").Value = 'Región de Arica y Parinacota'
$ws.Range("S16").Value = 1667
$ws.Range("T16").Value = 18

# Row 17
$ws.Range("D17").Value = 44432
$ws.Range("M17").Value = 10
$ws.Range("N17").Value = 35000
$ws.Range("O17").Value = 35000
$ws.Range("P17").Value = 35000
$ws.Range("Q17").Value = '$/caja 18 kilos'
$ws.Range("R17").Value = 'Perú'
$ws.Range("S17").Value = 1944
$ws.Range("T17").Value = 18

# Row 18
$ws.Range("D18").Value = 44431
$ws.Range("M18").Value = 30
$ws.Range("N18").Value = 35000
$ws.Range("O18").Value = 35000
$ws.Range("P18").Value = 35000
$ws.Range("Q18").Value = '$/caja 18 kilos'
$ws.Range("R18").Value = 'Región de Arica y Parinacota'
$ws.Range("S18").Value = 1944
$ws.Range("T18").Value = 18

# Row 19
$ws.Range("D19").Value = 44449
$ws.Range("M19").Value = 20
$ws.Range("N19").Value = 38000
$ws.Range("O19").Value = 38000
$ws.Range("P19").Value = 38000
$ws.Range("Q19").Value = '$/caja 18 kilos'
$ws.Range("R19").Value = 'Región de Arica y Parinacota'
$ws.Range("S19").Value = 2111
$ws.Range("T19").Value = 18

# Row 20
$ws.Range("D20").Value = 44424
$ws.Range("M20").Value = 15
$ws.Range("N20").Value = 35000
$ws.Range("O20").Value = 35000
$ws.Range("P20").Value = 35000
$ws.Range("Q20").Value = '$/caja 18 kilos'
$ws.Range("R20").Value = 'Región de Arica y Parinacota'
$ws.Range("S20").Value = 1944
$ws.Range("T20").Value = 18

# Row 21
$ws.Range("D21").Value = 44418
$ws.Range("M21").Value = 30
$ws.Range("N21").Value = 35000
$ws.Range("O21").Value = 35000
$ws.Range("P21").Value = 35000
$ws.Range("Q21").Value = '$/caja 18 kilos'
$ws.Range("R21").Value = 'Región de Arica y Parinacota'
$ws.Range("S21").Value = 1944
$ws.Range("T21").Value = 18

# Row 22
$ws.Range("D22").Value = 44704
$ws.Range("M22").Value = 25
$ws.Range("N22").Value = 35000
$ws.Range("O22").Value = 35000
$ws.Range("P22").Value = 35000
$ws.Range("Q22").Value = '$/caja 18 kilos'
$ws.Range("R22").Value = 'Región de Arica y Parinacota'
$ws.Range("S22").Value = 1944
$ws.Range("T22").Value = 18

